$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ============ 1. New cell values ============
$ws.Range("C4").Value = "en cours"
$ws.Range("B5").Value = "Quand on enchaine 2 accords dont l'intervalle entre les basses est une quarte ou une quinte, on garde la note commune et on enchaine les autres au plus près"
$ws.Range("B6").Value = "Quand on enchaine 2 accords dont l’intervalle entre les basses est une tierce (montante ou descendante), on garde les 2 notes communes et on enchaine les autres parties au plus près"
$ws.Range("B7").Value = "Quand on enchaine 2 accords dont l’intervalle entre les basses est une seconde, toutes les voix vont par mouvement contraire à la basse pour éviter les quintes et octaves consécutives"
$ws.Range("B8").Value = "Il est interdit de faire 2  quintes ou octaves consécutives même par mouvement contraire"
$ws.Range("C9").Value = "en cours"

# ============ 2. B7 left-aligned / indented cell (distinct style) ============
$ws.Range("B7").HorizontalAlignment = -4131
$ws.Range("B7").VerticalAlignment = -4108
$ws.Range("B7").WrapText = $true
$ws.Range("B7").IndentLevel = 4

# ============ 3. B9 rich text + left-aligned / indented + distinct font cell ============
# Set alignment/indent FIRST so the font change below does not fork an intermediate/orphan style.
$ws.Range("B9").Value = " L’accord de quinte diminuée du second degré du mode mineur et l’accord de septième degré du mode majeur seront de préférence utilisés à 3 voix, les doublures de notes dans ces accords étant d’un effet peu heureux"
$ws.Range("B9").HorizontalAlignment = -4131
$ws.Range("B9").VerticalAlignment = -4108
$ws.Range("B9").WrapText = $true
$ws.Range("B9").IndentLevel = 4
$ws.Range("B9").Font.Name = "Calibri"
$ws.Range("B9").Font.Size = 12
$run1 = $ws.Range("B9").Characters(1, 1)
$run1.Font.Name = "Times New Roman"
$run1.Font.Size = 7
$run2 = $ws.Range("B9").Characters(2, 212)
$run2.Font.Name = "Calibri"
$run2.Font.Size = 12

# ============ 4. C2 highlighted fill cell (distinct style) ============
$ws.Range("C2").HorizontalAlignment = -4108
$ws.Range("C2").VerticalAlignment = -4108
$ws.Range("C2").WrapText = $true
$ws.Range("C2").Interior.ThemeColor = 10

# ============ 5. Wrap text for the remaining plain center/center cells ============
# (looping one cell at a time -- this engine does not apply formatting to every area
#  of a multi-area/union Range, only the first, so avoid comma ranges here)
foreach ($addr in @("A2","B2","B3","C3","A4","B4","C4","B5","B6","B8","C9")) {
    $r = $ws.Range($addr)
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
    $r.WrapText = $true
}

# ============ 6. Row heights ============
$ws.Rows(5).RowHeight = 104

# ============ 7. Selection ============
$ws.Range("C10").Select()

Write-Host "done"
